$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Delete row 24 (content merged away / removed entirely) ---
$ws.Rows(24).Delete()

# --- Row 13: Programa resumido: / Semestral ---
$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"
$ws.Range("A13").EntireRow.RowHeight = 60

# --- Row 14: Short syllabus: / Generic biotechnological... ---
$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("B14").Value = "Generic biotechnological process, equipment sterilization, media sterilization by heating, air sterilization by filtration, kinetics and stoichiometry of microbial growth and products formation."
$ws.Range("C14").Value = "Generic biotechnological process, equipment sterilization, media sterilization by heating, air sterilization by filtration, kinetics and stoichiometry of microbial growth and products formation."
$ws.Range("A14").EntireRow.RowHeight = 60

# --- Row 15: Programa: / 01/01/2018 ---
$ws.Range("A15").Value = "Programa:"
$ws.Range("B15").Value = "01/01/2018"
$ws.Range("C15").Value = "01/01/2018"
$ws.Range("A15").EntireRow.RowHeight = 120

# --- Row 16: unchanged (Syllabus: / 1.Generic biotechnological...) ---
$ws.Range("A16").EntireRow.RowHeight = 120

# --- Row 17: Avaliação: only, clear B/C, remove custom height ---
$ws.Range("A17").Value = "Avaliação:"
$ws.Range("B17:C17").Clear()
$ws.Range("A17").EntireRow.AutoFit()

# --- Row 18: Método: / 5840876 - Walter de Carvalho ---
$ws.Range("A18").Value = "Método:"
$ws.Range("B18").Value = "5840876 - Walter de Carvalho"
$ws.Range("C18").Value = "5840876 - Walter de Carvalho"
$ws.Range("A18").EntireRow.RowHeight = 60

# --- Row 19: Critério: / Os alunos serão avaliados... ---
$ws.Range("A19").Value = "Critério:"
$ws.Range("B19").Value = "Os alunos serão avaliados formalmente por duas provas escritas (P1 e P2), sendo a segunda prova (P2) com peso 2."
$ws.Range("C19").Value = "Os alunos serão avaliados formalmente por duas provas escritas (P1 e P2), sendo a segunda prova (P2) com peso 2."
$ws.Range("A19").EntireRow.RowHeight = 60

# --- Row 20: Norma de recuperação: / A nota final (NF)... ---
$ws.Range("A20").Value = "Norma de recuperação:"
$ws.Range("B20").Value = "A nota final (NF) será calculada como: NF=(P1+(P2×2))/3. Serão aprovados os alunos que obtiverem NF maior ou igual 5,0."
$ws.Range("C20").Value = "A nota final (NF) será calculada como: NF=(P1+(P2×2))/3. Serão aprovados os alunos que obtiverem NF maior ou igual 5,0."
$ws.Range("A20").EntireRow.RowHeight = 60

# --- Row 21: Bibliografia: / Será oferecido um programa... ---
$ws.Range("A21").Value = "Bibliografia:"
$ws.Range("B21").Value = "Será oferecido um programa de recuperação avaliado por uma prova escrita final (PR).`nA média de recuperação (MR) será calculada como: MR=(NF+PR)/2. Serão aprovados os alunos que obtiverem MR maior ou igual a 5,0."
$ws.Range("C21").Value = "Será oferecido um programa de recuperação avaliado por uma prova escrita final (PR).`nA média de recuperação (MR) será calculada como: MR=(NF+PR)/2. Serão aprovados os alunos que obtiverem MR maior ou igual a 5,0."
$ws.Range("A21").EntireRow.RowHeight = 120

# --- Row 22: Requisitos: only, clear B/C, remove custom height ---
$ws.Range("A22").Value = "Requisitos:"
$ws.Range("B22:C22").Clear()
$ws.Range("A22").EntireRow.AutoFit()

# --- Row 23: B/C only = LOT2028 requirement text, no A cell ---
$ws.Range("B23").Value = "LOT2028 -  Tecnologia de Processos Fermentativos  (Requisito fraco)`n"
$ws.Range("C23").Value = "LOT2028 -  Tecnologia de Processos Fermentativos  (Requisito fraco)`n"
$ws.Range("A23").EntireRow.RowHeight = 30

